$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some updated price values look like plain decimal numbers (e.g. "27.00",
# "18.50") and would otherwise be auto-converted by Excel into numeric values,
# silently dropping the trailing zeros / formatting that the source feed
# provides as text. Force those specific cells to Text format first, cell
# by cell (multi-area ranges only honor the first area for formatting), so
# the assigned strings are preserved exactly.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"

# Write the refreshed values per row: B (coin name), C (link), D (price) and
# E (1h volume change).

# Row 2
$ws.Range("D2").Value = "43.412.26"
$ws.Range("E2").Value = "  -1.33%  "

# Row 3
$ws.Range("D3").Value = "2.375.18"
$ws.Range("E3").Value = "  +4.68%  "

# Row 4
$ws.Range("E4").Value = "  -0.34%  "

# Row 5
$ws.Range("D5").Value = "235.78"
$ws.Range("E5").Value = "  +1.04%  "

# Row 6
$ws.Range("D6").Value = "0.646"
$ws.Range("E6").Value = "  -1.81%  "

# Row 7
$ws.Range("D7").Value = "71.36"
$ws.Range("E7").Value = "  +11.90%  "

# Row 8
$ws.Range("E8").Value = "  -0.18%  "

# Row 9
$ws.Range("D9").Value = "0.472"
$ws.Range("E9").Value = "  +4.05%  "

# Row 10
$ws.Range("D10").Value = "0.0978"
$ws.Range("E10").Value = "  -0.73%  "

# Row 11
$ws.Range("D11").Value = "56.85"
$ws.Range("E11").Value = "  -2.28%  "

# Row 12
$ws.Range("D12").Value = "27.00"
$ws.Range("E12").Value = "  +0.12%  "

# Row 13
$ws.Range("D13").Value = "2.728.40"
$ws.Range("E13").Value = "  +4.80%  "

# Row 14
$ws.Range("E14").Value = "  +0.05%  "

# Row 15
$ws.Range("D15").Value = "16.02"
$ws.Range("E15").Value = "  +1.88%  "

# Row 16
$ws.Range("E16").Value = "  +1.73%  "

# Row 17
$ws.Range("D17").Value = "0.853"
$ws.Range("E17").Value = "  +1.52%  "

# Row 18
$ws.Range("D18").Value = "2.373.15"
$ws.Range("E18").Value = "  +4.76%  "

# Row 19
$ws.Range("D19").Value = "43.413.88"
$ws.Range("E19").Value = "  -1.21%  "

# Row 20
$ws.Range("E20").Value = "  +1.41%  "

# Row 21
$ws.Range("D21").Value = "74.73"
$ws.Range("E21").Value = "  +0.61%  "

# Row 22
$ws.Range("D22").Value = "6.35"
$ws.Range("E22").Value = "  +3.30%  "

# Row 23
$ws.Range("D23").Value = "250.57"
$ws.Range("E23").Value = "  -0.19%  "

# Row 24
$ws.Range("E24").Value = "  -0.02%  "

# Row 25
$ws.Range("D25").Value = "3.77"
$ws.Range("E25").Value = "  +13.27%  "

# Row 27
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "22.88"
$ws.Range("E27").Value = "  +1.73%  "

# Row 28
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "10.01"
$ws.Range("E28").Value = "  +0.37%  "

# Row 29
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "2.20"
$ws.Range("E29").Value = "  -4.57%  "

# Row 30
$ws.Range("D30").Value = "174.16"
$ws.Range("E30").Value = "  +0.07%  "

# Row 31
$ws.Range("E31").Value = "  +4.96%  "

# Row 32
$ws.Range("E32").Value = "  -6.14%  "

# Row 33
$ws.Range("D33").Value = "0.128"
$ws.Range("E33").Value = "  -0.01%  "

# Row 34
$ws.Range("E34").Value = "  -0.62%  "

# Row 35
$ws.Range("D35").Value = "0.0691"
$ws.Range("E35").Value = "  +0.63%  "

# Row 36
$ws.Range("D36").Value = "5.09"
$ws.Range("E36").Value = "  +1.71%  "

# Row 37
$ws.Range("E37").Value = "  +6.85%  "

# Row 38
$ws.Range("D38").Value = "6.56"
$ws.Range("E38").Value = "  +1.17%  "

# Row 39
$ws.Range("E39").Value = "  -1.52%  "

# Row 40
$ws.Range("E40").Value = "  -0.04%  "

# Row 41
$ws.Range("E41").Value = "  -0.07%  "

# Row 42
$ws.Range("D42").Value = "8.91"
$ws.Range("E42").Value = "  +1.03%  "

# Row 43
$ws.Range("D43").Value = "18.50"
$ws.Range("E43").Value = "  +6.61%  "

# Row 44
$ws.Range("D44").Value = "1.18"
$ws.Range("E44").Value = "  +7.39%  "

# Row 45
$ws.Range("D45").Value = "100.34"
$ws.Range("E45").Value = "  +1.66%  "

# Row 46
$ws.Range("E46").Value = "  +2.88%  "

# Row 47
$ws.Range("E47").Value = "  +2.09%  "

# Row 48
$ws.Range("D48").Value = "0.0956"
$ws.Range("E48").Value = "  +0.16%  "

# Row 49
$ws.Range("E49").Value = "  -7.89%  "

# Row 50
$ws.Range("D50").Value = "1.444.90"
$ws.Range("E50").Value = "  -0.78%  "

# Row 51
$ws.Range("D51").Value = "2.600.92"
$ws.Range("E51").Value = "  +4.95%  "

# Clear the temporary Text number format applied above so the cells revert
# to the workbook default (unstyled) cell format, leaving only the cell
# value itself changed, cell by cell as above.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D48").Style = "Normal"
